$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G column values per diff
$ws.Range("F495").Value = 10530
$ws.Range("F496").Value = 8277
$ws.Range("F497").Value = 7849
$ws.Range("F498").Value = 9270
$ws.Range("F503").Value = 7880
$ws.Range("F516").Value = 9534
$ws.Range("F517").Value = 7844
$ws.Range("F519").Value = 8101
$ws.Range("F526").Value = 9011
$ws.Range("F527").Value = 11715
$ws.Range("F533").Value = 11919
$ws.Range("F537").Value = 13734
$ws.Range("F542").Value = 10435
$ws.Range("F547").Value = 14103
$ws.Range("F551").Value = 18184
$ws.Range("F552").Value = 15443
$ws.Range("F553").Value = 15547
$ws.Range("F558").Value = 24723
$ws.Range("F559").Value = 22509
$ws.Range("F560").Value = 6111
$ws.Range("F561").Value = 24353
$ws.Range("F562").Value = 27114
$ws.Range("F563").Value = 14104
$ws.Range("F565").Value = 29008
$ws.Range("F567").Value = 23506
$ws.Range("F568").Value = 23976
$ws.Range("F580").Value = 28879
$ws.Range("F581").Value = 27121
$ws.Range("F583").Value = 29385
$ws.Range("F584").Value = 13257
$ws.Range("F588").Value = 25385
$ws.Range("F589").Value = 25828
$ws.Range("G589").Value = 475
$ws.Range("F590").Value = 28995
$ws.Range("G590").Value = 577
$ws.Range("F593").Value = 36999
$ws.Range("F597").Value = 29546
$ws.Range("G597").Value = 955
$ws.Range("F599").Value = 16601
$ws.Range("F600").Value = 39979
$ws.Range("F605").Value = 14800
$ws.Range("F608").Value = 45774
$ws.Range("G608").Value = 2891
$ws.Range("F609").Value = 36357
$ws.Range("G609").Value = 2154
$ws.Range("F610").Value = 33718
$ws.Range("G610").Value = 1921
$ws.Range("F611").Value = 33981
$ws.Range("G611").Value = 2124
$ws.Range("F612").Value = 16273
$ws.Range("G612").Value = 1427
$ws.Range("F613").Value = 21449
$ws.Range("G613").Value = 1882
$ws.Range("F614").Value = 47213
$ws.Range("G614").Value = 3317
$ws.Range("F615").Value = 36291
$ws.Range("G615").Value = 2338
$ws.Range("F616").Value = 37398
$ws.Range("G616").Value = 2560

# Append new rows 617-620
$ws.Range("A617").Value = 44511
$ws.Range("B617").Value = 547052
$ws.Range("C617").Value = 23203
$ws.Range("D617").Value = 6843
$ws.Range("E617").Value = 13446
$ws.Range("F617").Value = 37719
$ws.Range("G617").Value = 2596

$ws.Range("A618").Value = 44512
$ws.Range("B618").Value = 554296
$ws.Range("C618").Value = 26314
$ws.Range("D618").Value = 7244
$ws.Range("E618").Value = 13485
$ws.Range("F618").Value = 33024
$ws.Range("G618").Value = 2346

$ws.Range("A619").Value = 44513
$ws.Range("B619").Value = 560100
$ws.Range("C619").Value = 17941
$ws.Range("D619").Value = 5804
$ws.Range("E619").Value = 13537
$ws.Range("F619").Value = 15686
$ws.Range("G619").Value = 1667

$ws.Range("A620").Value = 44514
$ws.Range("B620").Value = 563445
$ws.Range("C620").Value = 10999
$ws.Range("D620").Value = 3345
$ws.Range("E620").Value = 13598
$ws.Range("F620").Value = 16532
$ws.Range("G620").Value = 1586

Write-Output "done"
